# This script applies the following changes to the "Artfynd" sheet:
#  1. The entire contents of row 12 and row 14 are swapped (all columns,
#     including cells that become blank/empty or cease to exist).
#  2. Row 21, column AX (Observatörer) has its two names re-ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full set of columns used anywhere on the sheet (A .. AY).
$cols = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T",
    "U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL",
    "AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY"
)

function Get-CellValue($ws, $col, $row) {
    return $ws.Range("$col$row").Value()
}

# Writes $val (captured from the "source" row) into column $col of row $row,
# taking care to faithfully reproduce three distinct states:
#   - $null            -> cell must not exist at all (ClearContents)
#   - empty string ""  -> cell must exist but be blank
#   - otherwise        -> cell must hold that value (number/bool/string)
function Set-CellValue($ws, $col, $row, $val) {
    $addr = "$col$row"
    if ($val -eq $null) {
        $ws.Range($addr).ClearContents()
    }
    elseif ($val -is [string]) {
        # Prefix text with a leading apostrophe so Excel always stores it
        # verbatim as text, even when it looks like a number/date/time
        # (e.g. "2026-02-16", "10:33") or is empty. Then reset the style so
        # no visible quote-prefix formatting is left on the cell.
        $ws.Range($addr).Value = "'" + $val
        $ws.Range($addr).Style = "Normal"
    }
    else {
        $ws.Range($addr).Value = $val
    }
}

# Capture the full "before" state of both rows first, so writing one row
# doesn't clobber data still needed for the other.
$row12vals = @{}
$row14vals = @{}
foreach ($col in $cols) {
    $row12vals[$col] = Get-CellValue $ws $col 12
    $row14vals[$col] = Get-CellValue $ws $col 14
}

# Swap: row 12 gets what row 14 used to hold, and vice versa.
foreach ($col in $cols) {
    Set-CellValue $ws $col 12 $row14vals[$col]
}
foreach ($col in $cols) {
    Set-CellValue $ws $col 14 $row12vals[$col]
}

# Row 21: reorder the observer names in the "Observatörer" column.
$ws.Range("AX21").Value = "Anna-Lena Thommson, Lars-Erik Nilsson"
